$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsCodebook = $wb.Worksheets.Item("Codebook")

$wsData.Range("D1").Value = "ShoeSize"
$wsData.Range("E1").Value = "EyeColor"

$wsCodebook.Range("A5").Value = "ShoeSize"
$wsCodebook.Range("A6").Value = "EyeColor"

$wsData.Range("E2").Select()
$wsCodebook.Range("A7").Select()
$wsCodebook.Activate()
